# Applies the StructureDefinition-employee-job-title.xlsx update:
#  - bump Version / Date metadata
#  - add Publisher + Jurisdiction values, dropping the stray duplicate
#    "Contact" row that used to sit under Publisher
#  - give the root Extension row its real Short/Definition text instead
#    of the generic FHIR boilerplate

$wb = $excel.ActiveWorkbook

# ---- Metadata sheet -------------------------------------------------
$meta = $wb.Worksheets.Item("Metadata")

$meta.Range("B3").Value = "6.0.0"
$meta.Range("B8").Value = "2022-01-21T20:46:54+00:00"
$meta.Range("B9").Value = "Alvearie Team"

# Row 10 used to be a second "Contact" / "No display for ContactDetail"
# row; turn it into "Jurisdiction" / "United States of America" ...
$meta.Range("A10").Value = "Jurisdiction"
$meta.Range("B10").Value = "United States of America"

# ... and remove the now-duplicate "Contact" row that followed it,
# shifting everything below up by one (21 rows -> 20 rows).
$meta.Rows.Item(11).Delete()

# ---- Elements sheet ---------------------------------------------------
$elements = $wb.Worksheets.Item("Elements")

# The root "Extension" element's Short/Definition previously held the
# generic base-Extension boilerplate; replace with the real text for
# this extension.
$elements.Range("K2").Value = "Employee Job Title"
$elements.Range("L2").Value = "Job title of the employee"
